$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (label text stays sequential; C/D/E values change) ---

# Row 8: line7 (was extr1)
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: line8 (was extr2)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16

# Row 10: extr1 (was extr3)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11: extr2 (was extr4)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12: extr3 (was extr5)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13: extr4 (was extr6)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14: extr5 (was extr7)
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15: extr6 (was extr8)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Add new rows 16 and 17 ---

# Row 16: extr7
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17: extr8
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
